$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2: mark task as done (copy formatting from an existing "done" cell, e.g. B4)
$ws.Range("B2").Value = "Cделано"
$ws.Range("B4").Copy()
$ws.Range("B2").PasteSpecial(-4122) # xlPasteFormats

# Move the active selection to B2
$ws.Range("B2").Select()
